$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 6 (dataset now has one fewer row)
$ws.Rows.Item(6).Delete()

# Adjust column widths 7 -> 8 for the affected columns
$ws.Columns.Item(2).ColumnWidth = 7.17
$ws.Columns.Item(3).ColumnWidth = 7.17
$ws.Columns.Item(7).ColumnWidth = 7.17
$ws.Columns.Item(10).ColumnWidth = 7.17
$ws.Columns.Item(11).ColumnWidth = 7.17
$ws.Columns.Item(12).ColumnWidth = 7.17
$ws.Columns.Item(13).ColumnWidth = 7.17
$ws.Columns.Item(15).ColumnWidth = 7.17
$ws.Columns.Item(17).ColumnWidth = 7.17
$ws.Columns.Item(22).ColumnWidth = 7.17
$ws.Columns.Item(24).ColumnWidth = 7.17
$ws.Columns.Item(27).ColumnWidth = 7.17
$ws.Columns.Item(28).ColumnWidth = 7.17
$ws.Columns.Item(29).ColumnWidth = 7.17
$ws.Columns.Item(30).ColumnWidth = 7.17
$ws.Columns.Item(34).ColumnWidth = 7.17

# Update data rows 2-5 with the new dataset
$ws.Range("A2").Value = 45079.50694444445
$ws.Range("B2").Value = 24.021
$ws.Range("C2").Value = 16.687
$ws.Range("D2").Value = 4.255
$ws.Range("E2").Value = 50.696
$ws.Range("F2").Value = 41.946
$ws.Range("G2").Value = 18.904
$ws.Range("H2").Value = 63.181
$ws.Range("I2").Value = 29.086
$ws.Range("J2").Value = 12.432
$ws.Range("K2").Value = 19.179
$ws.Range("L2").Value = 19.925
$ws.Range("M2").Value = 20.923
$ws.Range("N2").Value = 6.036
$ws.Range("O2").Value = 18.798
$ws.Range("P2").Value = 26.508
$ws.Range("Q2").Value = 15.598
$ws.Range("R2").Value = 3.832
$ws.Range("S2").Value = 2.606
$ws.Range("T2").Value = 278.79
$ws.Range("U2").Value = 52.368
$ws.Range("V2").Value = 17.351
$ws.Range("W2").Value = 34.892
$ws.Range("X2").Value = 18.126
$ws.Range("Y2").Value = 2.393
$ws.Range("Z2").Value = 31.289
$ws.Range("AA2").Value = 15.326
$ws.Range("AB2").Value = 13.706
$ws.Range("AC2").Value = 16.029
$ws.Range("AD2").Value = 20.711
$ws.Range("AE2").Value = 3.641
$ws.Range("AF2").Value = 55.941
$ws.Range("AG2").Value = 9.720000000000001
$ws.Range("AH2").Value = 21.693
$ws.Range("A3").Value = 45079.51388888889
$ws.Range("B3").Value = 13.932
$ws.Range("C3").Value = 9.795
$ws.Range("D3").Value = 1.732
$ws.Range("E3").Value = 29.697
$ws.Range("F3").Value = 24.609
$ws.Range("G3").Value = 10.964
$ws.Range("H3").Value = 44.715
$ws.Range("I3").Value = 16.87
$ws.Range("J3").Value = 7.275
$ws.Range("K3").Value = 11.086
$ws.Range("L3").Value = 11.799
$ws.Range("M3").Value = 12.414
$ws.Range("N3").Value = 3.504
$ws.Range("O3").Value = 10.903
$ws.Range("P3").Value = 15.39
$ws.Range("Q3").Value = 9.298
$ws.Range("R3").Value = 1.617
$ws.Range("S3").Value = 0.993
$ws.Range("T3").Value = 158.644
$ws.Range("U3").Value = 30.579
$ws.Range("V3").Value = 10.064
$ws.Range("W3").Value = 20.303
$ws.Range("X3").Value = 10.823
$ws.Range("Y3").Value = 1.365
$ws.Range("Z3").Value = 21.107
$ws.Range("AA3").Value = 8.888999999999999
$ws.Range("AB3").Value = 8.048999999999999
$ws.Range("AC3").Value = 9.416
$ws.Range("AD3").Value = 12.357
$ws.Range("AE3").Value = 1.294
$ws.Range("AF3").Value = 40.443
$ws.Range("AG3").Value = 5.6
$ws.Range("AH3").Value = 12.582
$ws.Range("A4").Value = 45079.52083333334
$ws.Range("B4").Value = 14.413
$ws.Range("C4").Value = 10.382
$ws.Range("D4").Value = 1.256
$ws.Range("E4").Value = 30.955
$ws.Range("F4").Value = 25.646
$ws.Range("G4").Value = 11.342
$ws.Range("H4").Value = 44.606
$ws.Range("I4").Value = 17.452
$ws.Range("J4").Value = 7.632
$ws.Range("K4").Value = 11.53
$ws.Range("L4").Value = 12.377
$ws.Range("M4").Value = 13.03
$ws.Range("N4").Value = 3.623
$ws.Range("O4").Value = 11.279
$ws.Range("P4").Value = 15.97
$ws.Range("Q4").Value = 9.590999999999999
$ws.Range("R4").Value = 1.073
$ws.Range("S4").Value = 0.752
$ws.Range("T4").Value = 164.352
$ws.Range("U4").Value = 31.55
$ws.Range("V4").Value = 10.411
$ws.Range("W4").Value = 21.061
$ws.Range("X4").Value = 11.261
$ws.Range("Y4").Value = 1.42
$ws.Range("Z4").Value = 21.276
$ws.Range("AA4").Value = 9.196
$ws.Range("AB4").Value = 8.253
$ws.Range("AC4").Value = 9.676
$ws.Range("AD4").Value = 13.005
$ws.Range("AE4").Value = 0.784
$ws.Range("AF4").Value = 40.218
$ws.Range("AG4").Value = 5.823
$ws.Range("AH4").Value = 13.016
$ws.Range("A5").Value = 45079.52777777778
$ws.Range("B5").Value = 14.41
$ws.Range("C5").Value = 10.49
$ws.Range("D5").Value = 1.04
$ws.Range("E5").Value = 31.05
$ws.Range("F5").Value = 25.72
$ws.Range("G5").Value = 11.34
$ws.Range("H5").Value = 44.57
$ws.Range("I5").Value = 17.45
$ws.Range("J5").Value = 7.69
$ws.Range("K5").Value = 11.56
$ws.Range("L5").Value = 12.45
$ws.Range("M5").Value = 13.11
$ws.Range("N5").Value = 3.62
$ws.Range("O5").Value = 11.28
$ws.Range("P5").Value = 16.01
$ws.Range("Q5").Value = 9.56
$ws.Range("R5").Value = 0.82
$ws.Range("S5").Value = 0.65
$ws.Range("T5").Value = 164.35
$ws.Range("U5").Value = 31.55
$ws.Range("V5").Value = 10.41
$ws.Range("W5").Value = 21.13
$ws.Range("X5").Value = 11.28
$ws.Range("Y5").Value = 1.42
$ws.Range("Z5").Value = 21.37
$ws.Range("AA5").Value = 9.199999999999999
$ws.Range("AB5").Value = 8.210000000000001
$ws.Range("AC5").Value = 9.640000000000001
$ws.Range("AD5").Value = 13.09
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 40.27
$ws.Range("AG5").Value = 5.84
$ws.Range("AH5").Value = 13.02
